$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Materials Identified"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Materials Identified")

$ws1.Range("C12").Formula = "=2534-868"
$ws1.Range("C13").Value = 38
$ws1.Range("C14").Value = 115
$ws1.Range("C15").Value = 236
$ws1.Range("C16").Value = 479

$ws1.Columns.Item(2).ColumnWidth = 45.42578125

# ---------------------------------------------------------------------------
# Sheet 2: "Cells identified"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Cells identified")

$ws2.Range("C3").Value = 40950
$ws2.Range("D3").Value = 39643
$ws2.Range("E3").Value = 37713

$ws2.Range("C14").Value = 37713
$ws2.Range("D14").Value = 39642
$ws2.Range("E14").Value = 40950

$ws2.Range("C15").Formula = "=`$E`$17-C14"
$ws2.Range("D15").Formula = "=`$E`$17-D14"
$ws2.Range("E15").Formula = "=`$E`$17-E14"

$ws2.Range("C19").Formula = "=C14/`$E`$17"
$ws2.Range("D19").Formula = "=D14/`$E`$17"
$ws2.Range("E19").Formula = "=E14/`$E`$17"

$wb.Application.Calculate()
